$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A14").Value = "Dividend"
$ws.Range("B20").Select() | Out-Null
